# Updated Random Forest so that it used less trees
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("New Results")

# 1) Grab the current (yellow-highlight) formatting of row 17 - the
#    previous "latest" Random Forest result - and stamp it onto the new
#    row 20 before row 17 itself gets reformatted below.
$ws.Range("A17:E17").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows(20).RowHeight = $ws.Rows(17).RowHeight

# 2) Row 17 is no longer the newest result, so it loses its highlight -
#    match the "no fill" formatting used by the other Random Forest rows
#    (e.g. row 13).
$ws.Range("A13:E13").Copy()
$ws.Range("A17:E17").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# 3) Fill in the new row's data - Random Forest retrained with fewer trees
#    (NumLearningCycles 100 instead of 499).
$ws.Range("A20").Value = "Random Forest"
$ws.Range("B20").Value = "MinLeafSize, 1, Method, Bag, NumLearningCycles, 100, (columns: model, year, mileage, fuel type, MPG, engine size)."
$ws.Range("C20").Value = 1172.3900000000001
$ws.Range("D20").Value = 836.57
$ws.Range("E20").Value = 0.94

$ws.Range("A17").Select()
$ws.Range("E20").Select()
